$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = -10.24160000000001
    18 = -11.96789999999999
    20 = -11.43100000000001
    27 = -12.43729999999999
    35 = -11.5161
    69 = -11.1002
    76 = -12.48170000000001
    78 = -12.8842
    82 = -12.00990000000001
    83 = -13.5156
    93 = -10.21549999999999
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}

$wb.Save()
